# "Generate Report for Handoff"
# Update the localization-status report to reflect that the content is
# now Ready for handoff (instead of In Translation), refresh the
# associated generation timestamps, and widen the Status columns so the
# new, longer status text displays properly.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus          # zh-cn status
$overview.Range("F2").Value = $newStatus          # de-de status
$overview.Range("G2").Value = "2016-09-07 07:23:16"   # Latest HO Xliff Generate Date

# Widen the zh-cn / de-de status columns (E, F) to fit "Ready for handoff"
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus              # Status
$zhcn.Range("H2").Value = "2016-09-07 07:23:08"   # Latest Handoff Datetime

# Widen the Status column (C) to fit "Ready for handoff"
$zhcn.Columns.Item(3).ColumnWidth = 16.33

# --- de-de sheet ------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus              # Status
$dede.Range("H2").Value = "2016-09-07 07:23:16"   # Latest Handoff Datetime

# Widen the Status column (C) to fit "Ready for handoff"
$dede.Columns.Item(3).ColumnWidth = 16.33
